$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Bruselas (repollito)" at the
# "Vega Central Mapocho de Santiago" market. It belongs chronologically right
# after the existing row 84, so insert a fresh row at position 85 - this shifts
# every subsequent row (old 85..111) down by one (new 86..112), matching the
# diff exactly, and Excel auto-extends the sheet dimension to A1:R112.
$ws.Rows.Item(85).Insert()

# Fill in the new row with the new observation's data.
$ws.Range("A85").Value = 9
$ws.Range("B85").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C85").Value = "Metropolitana"
$ws.Range("D85").Value = 45146
$ws.Range("E85").Value = 13
$ws.Range("F85").Value = 100112035
$ws.Range("G85").Value = "Bruselas (repollito)"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 52
$ws.Range("K85").Value = 17000
$ws.Range("L85").Value = 18000
$ws.Range("M85").Value = 17500
$ws.Range("N85").Value = "$/malla 15 kilos"
$ws.Range("O85").Value = "Provincia de Quillota"
$ws.Range("P85").Value = 1167
$ws.Range("Q85").Value = 15
$ws.Range("R85").Value = "Hortaliza"
